$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Rent/400/45858.22928240741 -> Trip/105/45889.22928240741 ---
$ws.Range("A2").Value = "Trip"
$ws.Range("B2").Value = 105
$ws.Range("C2").Value = 45889.22928240741

# --- Row 3: pizza/200/45858.22928240741 -> gh/890/45883.22928240741 ---
$ws.Range("A3").Value = "gh"
$ws.Range("B3").Value = 890
$ws.Range("C3").Value = 45883.22928240741

# --- Row 4 (new): asc/670/45881.22928240741 ---
$ws.Range("A4").Value = "asc"
$ws.Range("B4").Value = 670

# Copy the date cell formatting (style) from C2 onto C4, then set its value,
# so the new date cell matches the existing date-formatted column.
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = 45881.22928240741

$excel.CutCopyMode = $false
